$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for rows 3, 5, 8, 9
$ws.Range("F3").Value = 0
$ws.Range("F5").Value = -7
$ws.Range("F8").Value = 0
$ws.Range("F9").Value = 0
